$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.868554723999523
$ws.Range("C2").Value = 0.2230895546694285
$ws.Range("D2").Value = 0.02541487974600187
$ws.Range("E2").Value = 0.109308266228282
$ws.Range("F2").Value = 0.7337237621416648
$ws.Range("I2").Value = 0.6503881080935194
$ws.Range("L2").Value = 0.2035184872961509
$ws.Range("M2").Value = 0.2039161632909838
$ws.Range("N2").Value = 1.278446630484588
$ws.Range("O2").Value = 2.510429043626402
$ws.Range("B3").Value = 0.7866008231924866
$ws.Range("C3").Value = 0.2096252353765067
$ws.Range("D3").Value = 0.02394862812467125
$ws.Range("E3").Value = 0.110214494987166
$ws.Range("F3").Value = 0.7287116675615692
$ws.Range("I3").Value = 0.6558186271010022
$ws.Range("L3").Value = 0.200750244957085
$ws.Range("M3").Value = 0.1906489805552596
$ws.Range("N3").Value = 1.289411023549775
$ws.Range("O3").Value = 2.507811326002042
$ws.Range("B4").Value = 0.7363698976323292
$ws.Range("C4").Value = 0.2012869645223532
$ws.Range("D4").Value = 0.0230406219891961
$ws.Range("E4").Value = 0.1108039141966617
$ws.Range("F4").Value = 0.7261286202203507
$ws.Range("I4").Value = 0.6595363778698058
$ws.Range("L4").Value = 0.199151163346535
$ws.Range("M4").Value = 0.1825620133722765
$ws.Range("N4").Value = 1.296623658830093
$ws.Range("O4").Value = 2.507846823385449
$ws.Range("B5").Value = 0.7159241101786051
$ws.Range("C5").Value = 0.1978713240490748
$ws.Range("D5").Value = 0.02266868007409073
$ws.Range("E5").Value = 0.1110524143337162
$ws.Range("F5").Value = 0.7252002856644708
$ws.Range("I5").Value = 0.661147744388888
$ws.Range("L5").Value = 0.1985248895985805
$ws.Range("M5").Value = 0.1792815925659355
$ws.Range("N5").Value = 1.299683767271183
$ws.Range("O5").Value = 2.508274158580292
$ws.Range("B6").Value = 0.7125305743480794
$ws.Range("C6").Value = 0.1973030934987179
$ws.Range("D6").Value = 0.02260680389476022
$ws.Range("E6").Value = 0.1110941796437026
$ws.Range("F6").Value = 0.7250536412816899
$ws.Range("I6").Value = 0.6614211284806792
$ws.Range("L6").Value = 0.1984224310903286
$ws.Range("M6").Value = 0.178737798518398
$ws.Range("N6").Value = 1.300199200698508
$ws.Range("O6").Value = 2.508370047095752
$ws.Range("B7").Value = 0.7360940607948976
$ws.Range("C7").Value = 0.2012409715687369
$ws.Range("D7").Value = 0.02303561360681528
$ws.Range("E7").Value = 0.1108072319066045
$ws.Range("F7").Value = 0.7261155972353777
$ws.Range("I7").Value = 0.6595577192648321
$ws.Range("L7").Value = 0.1991426144173829
$ws.Range("M7").Value = 0.182517711112304
$ws.Range("N7").Value = 1.296664438882374
$ws.Range("O7").Value = 2.507850915182161
$ws.Range("B8").Value = 0.8402792903085015
$ws.Range("C8").Value = 0.2184619098920848
$ws.Range("D8").Value = 0.02491093020729096
$ws.Range("E8").Value = 0.1096138933733015
$ws.Range("F8").Value = 0.7318929501761033
$ws.Range("I8").Value = 0.652180945883778
$ws.Range("L8").Value = 0.2025431465323422
$ws.Range("M8").Value = 0.1993294885757209
$ws.Range("N8").Value = 1.282127486823967
$ws.Range("O8").Value = 2.509185352070716
$ws.Range("B9").Value = 1.045245462869218
$ws.Range("C9").Value = 0.2516623280739623
$ws.Range("D9").Value = 0.02852644670915794
$ws.Range("E9").Value = 0.1075350074939019
$ws.Range("F9").Value = 0.7471488321144477
$ws.Range("I9").Value = 0.6407597876439617
$ws.Range("L9").Value = 0.2100079451872432
$ws.Range("M9").Value = 0.232758920931019
$ws.Range("N9").Value = 1.257428932265874
$ws.Range("O9").Value = 2.524851093627717
$ws.Range("B10").Value = 1.196186840197299
$ws.Range("C10").Value = 0.2757021534020225
$ws.Range("D10").Value = 0.03114427024932809
$ws.Range("E10").Value = 0.1061661635196636
$ws.Range("F10").Value = 0.7607586634862429
$ws.Range("I10").Value = 0.6342287566908027
$ws.Range("L10").Value = 0.2159758550654374
$ws.Range("M10").Value = 0.2575932380129018
$ws.Range("N10").Value = 1.241599453196422
$ws.Range("O10").Value = 2.544341185788795
$ws.Range("B11").Value = 1.264921059430606
$ws.Range("C11").Value = 0.286560978408005
$ws.Range("D11").Value = 0.03232668827934049
$ws.Range("E11").Value = 0.1055777004257195
$ws.Range("F11").Value = 0.7674733401018301
$ws.Range("I11").Value = 0.631662408623292
$ws.Range("L11").Value = 0.2187954406617791
$ws.Range("M11").Value = 0.2689489737203345
$ws.Range("N11").Value = 1.23490022420436
$ws.Range("O11").Value = 2.554946461127315
$ws.Range("B12").Value = 1.290957879754217
$ws.Range("C12").Value = 0.2906617261313045
$ws.Range("D12").Value = 0.03277320792176397
$ws.Range("E12").Value = 0.1053597756810354
$ws.Range("F12").Value = 0.7700913820657149
$ws.Range("I12").Value = 0.630748857932744
$ws.Range("L12").Value = 0.2198781590974477
$ws.Range("M12").Value = 0.2732573227626105
$ws.Range("N12").Value = 1.232435478975688
$ws.Range("O12").Value = 2.559212844790693
$ws.Range("B13").Value = 1.285350022180467
$ws.Range("C13").Value = 0.2897790601004999
$ws.Range("D13").Value = 0.03267709727235513
$ws.Range("E13").Value = 0.1054064912877699
$ws.Range("F13").Value = 0.7695241881705073
$ws.Range("I13").Value = 0.630943014569695
$ws.Range("L13").Value = 0.2196443100371823
$ws.Range("M13").Value = 0.2723290826513747
$ws.Range("N13").Value = 1.232963099458026
$ws.Range("O13").Value = 2.558282861389756
$ws.Range("B14").Value = 1.2670629589025
$ws.Range("C14").Value = 0.2868985759581335
$ws.Range("D14").Value = 0.03236344862622076
$ws.Range("E14").Value = 0.1055596732035537
$ws.Range("F14").Value = 0.7676872177237328
$ws.Range("I14").Value = 0.6315860819830732
$ws.Range("L14").Value = 0.2188842162344287
$ws.Range("M14").Value = 0.2693032616048328
$ws.Range("N14").Value = 1.234696003169667
$ws.Range("O14").Value = 2.555292439001533
$ws.Range("B15").Value = 1.255862690243362
$ws.Range("C15").Value = 0.2851327253683564
$ws.Range("D15").Value = 0.0321711680962764
$ws.Range("E15").Value = 0.1056541411641403
$ws.Range("F15").Value = 0.7665718335318275
$ws.Range("I15").Value = 0.6319875701771807
$ws.Range("L15").Value = 0.2184205888016493
$ws.Range("M15").Value = 0.2674509170577295
$ws.Range("N15").Value = 1.235766846310753
$ws.Range("O15").Value = 2.553493337063315
$ws.Range("B16").Value = 1.191696187753791
$ws.Range("C16").Value = 0.2749909386332092
$ws.Range("D16").Value = 0.03106682464795796
$ws.Range("E16").Value = 0.1062053090173367
$ws.Range("F16").Value = 0.7603303829543506
$ws.Range("I16").Value = 0.6344046231387281
$ws.Range("L16").Value = 0.2157936912236238
$ws.Range("M16").Value = 0.2568522710005468
$ws.Range("N16").Value = 1.242047355693195
$ws.Range("O16").Value = 2.54368312856522
$ws.Range("B17").Value = 1.152349082234025
$ws.Range("C17").Value = 0.2687494286350045
$ws.Range("D17").Value = 0.03038716738265634
$ws.Range("E17").Value = 0.1065521936406091
$ws.Range("F17").Value = 0.7566355793734374
$ws.Range("I17").Value = 0.6359911055545311
$ws.Range("L17").Value = 0.2142089612491276
$ws.Range("M17").Value = 0.2503651587297782
$ws.Range("N17").Value = 1.246028711977097
$ws.Range("O17").Value = 2.538110539547205
$ws.Range("B18").Value = 1.129724363879916
$ws.Range("C18").Value = 0.2651522391388994
$ws.Range("D18").Value = 0.02999545313528529
$ws.Range("E18").Value = 0.1067549350036607
$ws.Range("F18").Value = 0.7545596964181556
$ws.Range("I18").Value = 0.6369416910655836
$ws.Range("L18").Value = 0.2133073300192905
$ws.Range("M18").Value = 0.2466394615732526
$ws.Range("N18").Value = 1.248365905097586
$ws.Range("O18").Value = 2.535069010940475
$ws.Range("B19").Value = 1.122065216821341
$ws.Range("C19").Value = 0.2639330549136218
$ws.Range("D19").Value = 0.02986268984016505
$ws.Range("E19").Value = 0.1068241334857567
$ws.Range("F19").Value = 0.7538652984473799
$ws.Range("I19").Value = 0.6372700812599987
$ws.Range("L19").Value = 0.2130037491013752
$ws.Range("M19").Value = 0.2453789593374651
$ws.Range("N19").Value = 1.249165350095545
$ws.Range("O19").Value = 2.534067303822553
$ws.Range("B20").Value = 1.156536966175679
$ws.Range("C20").Value = 0.2694145984398801
$ws.Range("D20").Value = 0.03045960033050932
$ws.Range("E20").Value = 0.1065149337281905
$ws.Range("F20").Value = 0.7570237977754317
$ws.Range("I20").Value = 0.6358182797765366
$ws.Range("L20").Value = 0.2143766380621912
$ws.Range("M20").Value = 0.2510551534617491
$ws.Range("N20").Value = 1.245600002728906
$ws.Range("O20").Value = 2.538686808921
$ws.Range("B21").Value = 1.272434088987893
$ws.Range("C21").Value = 0.2877449508236225
$ws.Range("D21").Value = 0.03245560857021701
$ws.Range("E21").Value = 0.1055145466824579
$ws.Range("F21").Value = 0.7682247352931455
$ws.Range("I21").Value = 0.6313956154276283
$ws.Range("L21").Value = 0.2191070675295066
$ws.Range("M21").Value = 0.2701917984747695
$ws.Range("N21").Value = 1.234185050617974
$ws.Range("O21").Value = 2.556164000409353
$ws.Range("B22").Value = 1.348229389395158
$ws.Range("C22").Value = 0.2996592177676121
$ws.Range("D22").Value = 0.03375289826920636
$ws.Range("E22").Value = 0.1048893712229622
$ws.Range("F22").Value = 0.7759843601248235
$ws.Range("I22").Value = 0.6288448261447783
$ws.Range("L22").Value = 0.222286094750217
$ws.Range("M22").Value = 0.282746264992241
$ws.Range("N22").Value = 1.227145003812211
$ws.Range("O22").Value = 2.56904599706175
$ws.Range("B23").Value = 1.307771925502948
$ws.Range("C23").Value = 0.2933064180544136
$ws.Range("D23").Value = 0.03306117807768771
$ws.Range("E23").Value = 0.1052204219239301
$ws.Range("F23").Value = 0.7718026987781172
$ws.Range("I23").Value = 0.6301751230646175
$ws.Range("L23").Value = 0.2205814094639464
$ws.Range("M23").Value = 0.2760414342385431
$ws.Range("N23").Value = 1.230863963343104
$ws.Range("O23").Value = 2.562036973912967
$ws.Range("B24").Value = 1.154643634773265
$ws.Range("C24").Value = 0.2691139027631664
$ws.Range("D24").Value = 0.0304268564146426
$ws.Range("E24").Value = 0.1065317686150977
$ws.Range("F24").Value = 0.756848133761693
$ws.Range("I24").Value = 0.6358962944045921
$ws.Range("L24").Value = 0.2143008019317705
$ws.Range("M24").Value = 0.2507431948957617
$ws.Range("N24").Value = 1.245793671854045
$ws.Range("O24").Value = 2.538425772231363
$ws.Range("B25").Value = 0.9897311952613563
$ws.Range("C25").Value = 0.2427422693873211
$ws.Range("D25").Value = 0.02755506235203597
$ws.Range("E25").Value = 0.1080695106551994
$ws.Range("F25").Value = 0.7426005586038187
$ws.Range("I25").Value = 0.6435231486489705
$ws.Range("L25").Value = 0.2079034657320875
$ws.Range("M25").Value = 0.2236667525887981
$ws.Range("N25").Value = 1.263703380562539
$ws.Range("O25").Value = 2.519213474105385
